# Update the attendance summary sheet: mark "Total Attendance Count" and
# "Real" (and in special cases "Invalid"/"Absent") columns from 0 to 1
# for each date row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where Total Attendance Count (D) and Real (E) become 1
$rowsDE = @(4, 5, 6, 9, 10, 11, 12, 13, 15, 17, 18)
foreach ($r in $rowsDE) {
    $ws.Range("D$r").Value = 1
    $ws.Range("E$r").Value = 1
}

# Rows where only Absent (H) becomes 1
$rowsH = @(7, 8, 14, 16)
foreach ($r in $rowsH) {
    $ws.Range("H$r").Value = 1
}

# Row 3 is special: Invalid (G) and Absent (H) become 1
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1
